# Update the build-version timestamp that is stamped throughout the workbook
# from "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST".

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- "About" sheet ---------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  " + '"Global Energy Monitor, Coal mine boundaries and methane sources for Jinjiazhuang Coal Mine, China, M1999, version ' + "'$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet --------------------------------
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 8; $row++) {
    $wsData.Cells.Item($row, 19).Value = $newVersion
}
